# "default to allow new next piece to be queued and use auto_ as default"
#
# The "cue" column (K) was previously only populated on the header row.
# Mirror the "stage" value already present in column I (the stage/next
# piece tag, e.g. "1b", "1c") into the new "cue" cells K2/K3 so a next
# piece/cue is queued by default.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = $ws.Range("I2").Value2
$ws.Range("K3").Value = $ws.Range("I3").Value2

# Move/save the active selection onto the newly populated cue cell.
$ws.Range("K3").Select() | Out-Null
